# edit.ps1
# Applies the changes described by the diff to the "Main" worksheet:
#  1. Font "Helvetica" -> "Times New Roman" for the two fonts used in the sheet
#     (the bold white header font, and the plain body font).
#  2. Re-apply the two solid fill colors used for row striping (cream / light
#     gray) so they get re-serialized (values are unchanged, 0xFFFEC8 and
#     0xEEEEEE respectively).
#  3. Column J width 24 -> ~16.8 characters.
#  4. Column J ("Preference") data values replaced with a new set of numbers
#     for rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fonts -------------------------------------------------------------
# Font id 2 (bold, white, used by the header row A1:J1) and font id 3 (plain,
# used by every data cell A2:J26) both switch from Helvetica to Times New
# Roman. Applying to the full used range covers every cell that references
# either font.
$ws.Range("A1:J26").Font.Name = "Times New Roman"

# --- 2. Fill colors ---------------------------------------------------------
# Rows 2-6 use the cream fill (RGB FFFEC8); starting at row 8 the even rows
# use the light-gray fill (RGB EEEEEE). Re-assert these colors (same value)
# so the fills get rewritten.
$ws.Range("A2:J6").Interior.Color = 13172479   # RGB(0xFE,0xFE,0xC8) -> BGR packed, matches fgColor 00FFFEC8
$ws.Range("A8:J8,A10:J10,A12:J12,A14:J14,A16:J16,A18:J18,A20:J20,A22:J22,A24:J24,A26:J26").Interior.Color = 15658734  # RGB(0xEE,0xEE,0xEE)

# --- 3. Column width --------------------------------------------------------
# Column J ("Preference") narrows from 24 to 16.8 characters. The engine
# quantizes ColumnWidth to whole pixels (using the workbook's Normal font
# metrics) before storing it, so an input of 16.8 itself would actually be
# stored as 17.67; requesting 16 is what rounds to the closest achievable
# stored width to 16.8 (16.83).
$ws.Columns.Item(10).ColumnWidth = 16

# --- 4. Column J values ------------------------------------------------------
$ws.Range("J2").Value = 0.109485
$ws.Range("J3").Value = 0.765132
$ws.Range("J4").Value = 0.304068
$ws.Range("J5").Value = 0.523421
$ws.Range("J6").Value = 0.396557
$ws.Range("J7").Value = 0.273592
$ws.Range("J8").Value = 0.268625
$ws.Range("J9").Value = 0.017357
$ws.Range("J10").Value = 0.18808
$ws.Range("J11").Value = 0.889866
$ws.Range("J12").Value = 0.843898
$ws.Range("J13").Value = 0.309
$ws.Range("J14").Value = 0.292024
$ws.Range("J15").Value = 0.47706
$ws.Range("J16").Value = 0.878262
$ws.Range("J17").Value = 0.806481
$ws.Range("J18").Value = 0.230078
$ws.Range("J19").Value = 0.048044
$ws.Range("J20").Value = 0.106818
$ws.Range("J21").Value = 0.647976
$ws.Range("J22").Value = 0.6198900000000001
$ws.Range("J23").Value = 0.8619520000000001
$ws.Range("J24").Value = 0.919411
$ws.Range("J25").Value = 0.590913
$ws.Range("J26").Value = 0.183695
